# Update the "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets. Both sheets share the same row layout for
# these rows, so the same set of updates is applied to each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 206
    3  = 249
    4  = 267
    5  = 801
    7  = 6142
    8  = 44
    11 = 62
    14 = 186
    15 = 427
    16 = 36
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
